# Order upgrade - second try at reformatting - add options for order
#
# Adds a second worksheet ("Sheet2") with a small A/B table, formats its
# header row (bold, centered/top-aligned, boxed with a thin border), and
# leaves the original "Sheet" untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Create the new worksheet. Worksheets.Add() inserts a blank sheet before
# the active one, so populate it first (writes against the fresh reference
# work reliably) and only rename/reposition it afterwards.
$new = $wb.Worksheets.Add()

# --- Header row ---------------------------------------------------------
$new.Range("A1").Value = "A"
$new.Range("B1").Value = "B"

# --- Data rows -----------------------------------------------------------
$new.Range("A2").Value = 1
$new.Range("B2").Value = 4
$new.Range("A3").Value = 2
$new.Range("B3").Value = 5
$new.Range("A4").Value = 3
$new.Range("B4").Value = 6

# --- Header formatting: bold, centered + top aligned, thin box border ---
$header = $new.Range("A1:B1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1

# --- Name and position the sheet after "Sheet" ---------------------------
$new.Name = "Sheet2"
$new.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))
